$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- "Gate Entry Page added": drop the two placeholder sheets, keep only "Data" ---
$null = $wb.Worksheets.Item("Data1").Delete()
$null = $wb.Worksheets.Item("Data2").Delete()

$ws = $wb.Worksheets.Item("Data")

# --- Drop the hyperlinks that belonged to the old rows 3 & 4 (keep the row-2 ones) ---
$hls = @($ws.Hyperlinks)
for ($i = $hls.Count - 1; $i -ge 2; $i--) {
    $null = $hls[$i].Delete()
}

# --- Remove old rows 3 & 4 (TestCase002 / TestCase003 rows) ---
$null = $ws.Rows.Item(3).Delete()
$null = $ws.Rows.Item(3).Delete()

# --- Clear leftover formatting so every cell starts from the plain default style ---
$ws.Cells.ClearFormats()

# --- Rewrite header row (row 1) for the new Gate Entry payload columns ---
$ws.Range("A1").Value = "TestCases"
$ws.Range("B1").Value = "type"
$ws.Range("C1").Value = "username"
$ws.Range("D1").Value = "password"
$ws.Range("E1").Value = "docType"
$ws.Range("F1").Value = "documentNo"
$ws.Range("G1").Value = "invoiceNum"
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "vehicleNum"
$ws.Range("J1").Value = "name"
$ws.Range("K1").Value = "number"

# --- Rewrite data row (row 2) ---
$ws.Range("A2").Value = "TestCase001- Gate Entry with the Complete Qty receiving of all the material codes from the selected order."
$ws.Range("B2").Value = "Positive"
$ws.Range("C2").Value = "Nikhil10"
$ws.Range("D2").Value = "Bcil@12345678"
$ws.Range("E2").Value = "PO"
$ws.Range("F2").Value = 24032502
$ws.Range("G2").Value = "TNTINV8910"
$ws.Range("H2").Value = "24-052025"
$ws.Range("I2").Value = "DL8SQQ7313"
$ws.Range("J2").Value = "Driver"
$ws.Range("K2").Value = 9999999999

# --- Formatting: wrap text across the used area, taller rows for the wrapped text ---
$ws.Range("A1:K2").WrapText = $true
$ws.Rows.Item(1).RowHeight = 28.5
$ws.Rows.Item(2).RowHeight = 57

# --- Column widths for the new layout ---
$ws.Columns.Item(1).ColumnWidth = 24.33203125
$ws.Columns.Item(2).ColumnWidth = 7.19921875
$ws.Columns.Item(3).ColumnWidth = 8.1328125
$ws.Columns.Item(4).ColumnWidth = 7.46484375
$ws.Columns.Item(5).ColumnWidth = 19.46484375
$ws.Columns.Item(8).ColumnWidth = 4.265625
$ws.Columns.Item(11).ColumnWidth = 10.73046875

# --- Selection matches the authored state ---
$ws.Range("J2:J5").Select()

Write-Host "Gate Entry Page added"
